$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of product data (row 6) - a "category without children" style entry
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Klangschalen"
$ws.Range("C6").Value = "Testschale"
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "KKS003.JPG"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("G6").Value = "Klangschale beschreibung"

# Update the active selection to reflect where the user ended up after editing
[void]$ws.Range("G7").Select()
